$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 1234, shifting existing rows 1234..1322 down to 1235..1323
$ws.Rows.Item(1234).Insert()

# Fill the new row 1234 with the new data record.
# Columns A, B, C, E, F, G, N, O, Q, R mirror the (unchanged) constants
# shared by every record in this dataset / this particular row's group.
$ws.Cells.Item(1234, 1).Value = 10
$ws.Cells.Item(1234, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1234, 3).Value = "La Araucanía"
$ws.Cells.Item(1234, 4).Value = 45265
$ws.Cells.Item(1234, 5).Value = 9
$ws.Cells.Item(1234, 6).Value = 100114001
$ws.Cells.Item(1234, 7).Value = "Papa"
$ws.Cells.Item(1234, 8).Value = "Cornado"
$ws.Cells.Item(1234, 9).Value = "1a nueva(o)"
$ws.Cells.Item(1234, 10).Value = 750
$ws.Cells.Item(1234, 11).Value = 22000
$ws.Cells.Item(1234, 12).Value = 25000
$ws.Cells.Item(1234, 13).Value = 23400
$ws.Cells.Item(1234, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(1234, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(1234, 16).Value = 936
$ws.Cells.Item(1234, 17).Value = 25
$ws.Cells.Item(1234, 18).Value = "Hortaliza"
